$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E, shifting existing columns E.. to the right.
$ws.Range("E1").EntireColumn.Insert()

# Populate the new column header and data.
$ws.Range("E1").Value = "pt_max"
$ws.Range("E2:E9").Value = 50

# Update selection to match the filled range (mirrors the saved UI state).
$ws.Range("E2:E9").Select()
